# Update the cryptos list with refreshed price/volume data (and a couple
# of rows whose rank swapped places, bringing new Coin/Link values along).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to stay a literal text string even when the value
    # looks numeric (e.g. "541.24"), mirroring how these price/volume
    # columns are stored as plain text in the source data.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Rows where only the Volume(1h) column (E) changed.
Set-TextValue $ws.Range("E4")  "  +0.17%  "
Set-TextValue $ws.Range("E7")  "  +0.13%  "
Set-TextValue $ws.Range("E26") "  -0.13%  "
Set-TextValue $ws.Range("E29") "  +0.19%  "

# Rows where Price (D) and Volume(1h) (E) changed.
$priceVolumeUpdates = @(
    @{ Row = 2;  D = "60.947.49";  E = "  -4.72%  " },
    @{ Row = 3;  D = "2.962.76";   E = "  -4.28%  " },
    @{ Row = 5;  D = "541.24";     E = "  -0.70%  " },
    @{ Row = 6;  D = "130.29";     E = "  -7.35%  " },
    @{ Row = 8;  D = "2.962.49";   E = "  -4.15%  " },
    @{ Row = 9;  D = "0.489";      E = "  -1.97%  " },
    @{ Row = 10; D = "0.144";      E = "  -7.86%  " },
    @{ Row = 11; D = "5.82";       E = "  -11.20%  " },
    @{ Row = 12; D = "0.440";      E = "  -4.07%  " },
    @{ Row = 13; D = "0.0000216";  E = "  -4.35%  " },
    @{ Row = 14; D = "33.47";      E = "  -4.17%  " },
    @{ Row = 15; D = "3.452.81";   E = "  -3.84%  " },
    @{ Row = 18; D = "2.974.73";   E = "  -3.79%  " },
    @{ Row = 19; D = "6.51";       E = "  -2.52%  " },
    @{ Row = 20; D = "461.64";     E = "  -3.97%  " },
    @{ Row = 21; D = "12.98";      E = "  -3.67%  " },
    @{ Row = 22; D = "0.659";      E = "  -6.09%  " },
    @{ Row = 23; D = "6.87";       E = "  -3.52%  " },
    @{ Row = 24; D = "79.18";      E = "  +0.12%  " },
    @{ Row = 25; D = "11.78";      E = "  -4.86%  " },
    @{ Row = 27; D = "2.67";       E = "  -2.16%  " },
    @{ Row = 28; D = "7.53";       E = "  -7.17%  " },
    @{ Row = 30; D = "1.86";       E = "  -2.69%  " },
    @{ Row = 31; D = "25.10";      E = "  -4.88%  " },
    @{ Row = 32; D = "1.11";       E = "  -4.38%  " },
    @{ Row = 33; D = "2.25";       E = "  -4.84%  " },
    @{ Row = 36; D = "5.76";       E = "  -4.56%  " },
    @{ Row = 37; D = "439.91";     E = "  -11.64%  " },
    @{ Row = 38; D = "3.131.23";   E = "  -3.86%  " },
    @{ Row = 39; D = "0.0779";     E = "  -2.97%  " },
    @{ Row = 40; D = "0.0373";     E = "  -7.91%  " },
    @{ Row = 41; D = "0.116";      E = "  -3.35%  " },
    @{ Row = 42; D = "7.98";       E = "  -1.95%  " },
    @{ Row = 45; D = "25.33";      E = "  +0.00%  " },
    @{ Row = 46; D = "0.237";      E = "  -6.75%  " },
    @{ Row = 47; D = "0.106";      E = "  -2.75%  " },
    @{ Row = 48; D = "116.00";     E = "  -6.66%  " },
    @{ Row = 49; D = "1.92";       E = "  -6.46%  " },
    @{ Row = 50; D = "1.28";       E = "  +7.51%  " },
    @{ Row = 51; D = "0.0₃0475";   E = "  -11.35%  " }
)

foreach ($u in $priceVolumeUpdates) {
    Set-TextValue $ws.Range("D$($u.Row)") $u.D
    $ws.Range("E$($u.Row)").Value = $u.E
}

# Rows 16 and 17 swapped ranking: row 16 now holds what was row 17's coin
# (WrappedBTC) and vice versa, each with its own freshly updated price and
# volume figures.
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "61.078.73"
$ws.Range("E16").Value = "  -4.57%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.109"
$ws.Range("E17").Value = "  -3.26%  "

# Rows 34 and 35 swapped ranking: OKB <-> NEARProtocol.
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "5.35"
$ws.Range("E34").Value = "  -1.02%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D35") "54.10"
$ws.Range("E35").Value = "  -5.50%  "

# Rows 43 and 44 swapped ranking: dogwifhat <-> USDe.
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D44") "2.34"
$ws.Range("E44").Value = "  -13.71%  "
